{"js": "// Strengthen the \"KEY ACHIEVEMENTS AND IMPACT\" section: turn the single\n// \"Software Development and Innovation\" bullet group into three headed\n// groups (\"Technical Innovation & Platform Development\",\n// \"Data Engineering & Analytics\", \"Research Leadership & Client Success\")\n// with many more bullet points, per the diff.\n\nconst body = context.document.body;\n\n// --- Locate the anchor paragraphs we need to touch -------------------\nconst heading3Results = body.search(\"Software Development and Innovation\", {\n  matchCase: true,\n  matchWholeWord: false,\n});\nconst bullet1Results = body.search(\n  \"\\u2022 Conceived and deployed redistricting software used by thousands of analysts nationwide\",\n  { matchCase: true, matchWholeWord: false }\n);\nconst econSimResults = body.search(\n  \"\\u2022 Created econometric simulation platform for humanitarian intervention modeling\",\n  { matchCase: true, matchWholeWord: false }\n);\nconst surveyOpsResults = body.search(\n  \"\\u2022 Built comprehensive survey operations platform from RFP through deployment\",\n  { matchCase: true, matchWholeWord: false }\n);\n\nheading3Results.load(\"items\");\nbullet1Results.load(\"items\");\neconSimResults.load(\"items\");\nsurveyOpsResults.load(\"items\");\nawait context.sync();\n\nconst heading3Para = heading3Results.items[0].paragraphs.getFirst();\nconst bullet1Para = bullet1Results.items[0].paragraphs.getFirst();\nconst econSimPara = econSimResults.items[0].paragraphs.getFirst();\nconst surveyOpsPara = surveyOpsResults.items[0].paragraphs.getFirst();\nawait context.sync();\n\n// --- 1) Rename the section heading ------------------------------------\nheading3Para.insertText(\n  \"Technical Innovation & Platform Development\",\n  Word.InsertLocation.replace\n);\n\n// --- 2) Strengthen the first bullet and add four new ones after it ----\nbullet1Para.insertText(\n  \"\\u2022 Conceived, architected, engineered and deployed cloud-based redistricting software used by thousands of analysts nationwide\",\n  Word.InsertLocation.replace\n);\n\nlet afterBullet1 = bullet1Para.insertParagraph(\n  \"\\u2022 Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party\",\n  Word.InsertLocation.after\n);\nafterBullet1 = afterBullet1.insertParagraph(\n  \"\\u2022 Developed RACSO platform for pollsters to fully administer research, analyzing bids from 1,200 vendors\",\n  Word.InsertLocation.after\n);\nafterBullet1 = afterBullet1.insertParagraph(\n  \"\\u2022 Engineered FLEEM system using Twilio API for thousands of simultaneous phone calls for IVR polls\",\n  Word.InsertLocation.after\n);\nafterBullet1 = afterBullet1.insertParagraph(\n  \"\\u2022 Created SimCrisis platform for humanitarian intervention modeling used by International Red Cross and UNICEF\",\n  Word.InsertLocation.after\n);\n\n// (The \"Developed boundary estimation system...\" bullet right after stays\n// exactly as-is, now trailing the new Technical Innovation group.)\n\n// --- 3) Turn the old \"Created econometric simulation...\" bullet into the\n//        new \"Data Engineering & Analytics\" heading, then add its bullets.\n//        Insert the new bullet paragraphs *before* changing the anchor's\n//        style, otherwise they inherit the Heading 3 style too.\nlet afterDataEng = econSimPara.insertParagraph(\n  \"\\u2022 Designed, architected and created multi-tenant data warehouse tracking decades of political, geographical, econometric change\",\n  Word.InsertLocation.after\n);\nafterDataEng = afterDataEng.insertParagraph(\n  \"\\u2022 Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%\",\n  Word.InsertLocation.after\n);\nafterDataEng = afterDataEng.insertParagraph(\n  \"\\u2022 Developed advanced data pipelines for machine learning applications enhancing consumer segmentation and predictive modeling\",\n  Word.InsertLocation.after\n);\nafterDataEng = afterDataEng.insertParagraph(\n  \"\\u2022 Built fraud detection systems for campaign finance data analysis across multi-terabyte datasets\",\n  Word.InsertLocation.after\n);\nafterDataEng = afterDataEng.insertParagraph(\n  \"\\u2022 Transformed small data team into big data engineering team using Hadoop Clusters and Hive on AWS\",\n  Word.InsertLocation.after\n);\nafterDataEng = afterDataEng.insertParagraph(\n  \"\\u2022 Introduced version control and Agile methodologies, improving project delivery timelines by 40%\",\n  Word.InsertLocation.after\n);\n\neconSimPara.insertText(\"Data Engineering & Analytics\", Word.InsertLocation.replace);\neconSimPara.style = \"Heading 3\";\n\n// --- 4) Insert the new \"Research Leadership & Client Success\" heading and\n//        its first three bullets right before the untouched\n//        \"Built comprehensive survey operations...\" bullet. Again, add the\n//        bullet paragraphs first, then style the heading paragraph.\nlet researchHeading = surveyOpsPara.insertParagraph(\n  \"Research Leadership & Client Success\",\n  Word.InsertLocation.before\n);\n\nlet beforeSurveyOps = researchHeading.insertParagraph(\n  \"\\u2022 Led multi-million dollar research projects involving sensitive consumer data with privacy compliance\",\n  Word.InsertLocation.after\n);\nbeforeSurveyOps = beforeSurveyOps.insertParagraph(\n  \"\\u2022 Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders\",\n  Word.InsertLocation.after\n);\nbeforeSurveyOps = beforeSurveyOps.insertParagraph(\n  \"\\u2022 Delivered actionable consumer insights and market intelligence for political candidates and major organizations\",\n  Word.InsertLocation.after\n);\n\nresearchHeading.style = \"Heading 3\";\n\n// --- 5) Add the two trailing bullets after the untouched survey-ops bullet\nlet afterSurveyOps = surveyOpsPara.insertParagraph(\n  \"\\u2022 Regular expert testimony and source on public opinion for journalists, elected officials, and NGO leadership\",\n  Word.InsertLocation.after\n);\nafterSurveyOps = afterSurveyOps.insertParagraph(\n  \"\\u2022 Redistricting analysis used in court cases with rigorous methodology and expert testimony\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Strengthen the \"KEY ACHIEVEMENTS AND IMPACT\" section: turn the single\n# \"Software Development and Innovation\" bullet group into three headed\n# groups (\"Technical Innovation & Platform Development\",\n# \"Data Engineering & Analytics\", \"Research Leadership & Client Success\")\n# with many more bullet points, per the diff.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphByText($doc, $text) {\n    # NOTE: always call this right before using the result and re-call it\n    # after any document mutation \u2014 Paragraph/Range references captured\n    # earlier can go stale (drift to the wrong paragraph) once paragraphs\n    # are inserted/removed elsewhere in the document.\n    $range = $doc.Content\n    $range.Find.Execute($text) | Out-Null\n    return $range.Paragraphs(1)\n}\n\n# Helper: insert a new paragraph with $text right after $afterPara and\n# return the newly created paragraph (so callers can chain more inserts).\nfunction Add-ParagraphAfter($afterPara, $text) {\n    $afterPara.Range.InsertParagraphAfter()\n    $newPara = $afterPara.Next()\n    $newPara.Range.Text = $text\n    return $newPara\n}\n\n# --- 1) Rename the section heading --------------------------------------\n$heading3Para = Find-ParagraphByText $d \"Software Development and Innovation\"\n$heading3Para.Range.Text = \"Technical Innovation & Platform Development\"\n\n# --- 2) Strengthen the first bullet and add four new ones after it -----\n$bullet1Para = Find-ParagraphByText $d \"Conceived and deployed redistricting software used by thousands of analysts nationwide\"\n$bullet1Para.Range.Text = \"\u2022 Conceived, architected, engineered and deployed cloud-based redistricting software used by thousands of analysts nationwide\"\n\n$cursor = Add-ParagraphAfter $bullet1Para \"\u2022 Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party\"\n$cursor = Add-ParagraphAfter $cursor \"\u2022 Developed RACSO platform for pollsters to fully administer research, analyzing bids from 1,200 vendors\"\n$cursor = Add-ParagraphAfter $cursor \"\u2022 Engineered FLEEM system using Twilio API for thousands of simultaneous phone calls for IVR polls\"\n$cursor = Add-ParagraphAfter $cursor \"\u2022 Created SimCrisis platform for humanitarian intervention modeling used by International Red Cross and UNICEF\"\n\n# (The \"Developed boundary estimation system...\" bullet right after stays\n# exactly as-is, now trailing the new Technical Innovation group.)\n\n# --- 3) Turn the old \"Created econometric simulation...\" bullet into the\n#        new \"Data Engineering & Analytics\" heading, then add its bullets.\n#        Re-find this anchor now since earlier inserts shifted the document.\n$econSimPara = Find-ParagraphByText $d \"Created econometric simulation platform for humanitarian intervention modeling\"\n\n$cursor = Add-ParagraphAfter $econSimPara \"\u2022 Designed, architected and created multi-tenant data warehouse tracking decades of political, geographical, econometric change\"\n$cursor = Add-ParagraphAfter $cursor \"\u2022 Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%\"\n$cursor = Add-ParagraphAfter $cursor \"\u2022 Developed advanced data pipelines for machine learning applications enhancing consumer segmentation and predictive modeling\"\n$cursor = Add-ParagraphAfter $cursor \"\u2022 Built fraud detection systems for campaign finance data analysis across multi-terabyte datasets\"\n$cursor = Add-ParagraphAfter $cursor \"\u2022 Transformed small data team into big data engineering team using Hadoop Clusters and Hive on AWS\"\n$cursor = Add-ParagraphAfter $cursor \"\u2022 Introduced version control and Agile methodologies, improving project delivery timelines by 40%\"\n\n$econSimPara.Range.Text = \"Data Engineering & Analytics\"\n$econSimPara.Style = \"Heading 3\"\n\n# --- 4) Insert the new \"Research Leadership & Client Success\" heading and\n#        its first three bullets right before the untouched\n#        \"Built comprehensive survey operations...\" bullet. Re-find this\n#        anchor now since earlier inserts shifted the document. After\n#        InsertParagraphBefore(), the anchor paragraph object *itself*\n#        becomes the new (empty) paragraph, and the original content\n#        shifts to .Next() \u2014 so write the new heading into $surveyOpsPara\n#        directly and use its .Next() as the real survey-ops paragraph.\n$surveyOpsPara = Find-ParagraphByText $d \"Built comprehensive survey operations platform from RFP through deployment\"\n$surveyOpsPara.Range.InsertParagraphBefore()\n$researchHeading = $surveyOpsPara\n$researchHeading.Range.Text = \"Research Leadership & Client Success\"\n\n$cursor = Add-ParagraphAfter $researchHeading \"\u2022 Led multi-million dollar research projects involving sensitive consumer data with privacy compliance\"\n$cursor = Add-ParagraphAfter $cursor \"\u2022 Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders\"\n$cursor = Add-ParagraphAfter $cursor \"\u2022 Delivered actionable consumer insights and market intelligence for political candidates and major organizations\"\n\n# Style the heading *after* its children exist \u2014 Add-ParagraphAfter()\n# clones the preceding paragraph's formatting for the new paragraph mark,\n# so styling the anchor first would make the new bullets inherit Heading 3.\n$researchHeading.Style = \"Heading 3\"\n\n# --- 5) Add the two trailing bullets after the untouched survey-ops bullet.\n#        Re-find this anchor one more time since step 4 inserted a\n#        paragraph right before it.\n$surveyOpsPara = Find-ParagraphByText $d \"Built comprehensive survey operations platform from RFP through deployment\"\n$cursor = Add-ParagraphAfter $surveyOpsPara \"\u2022 Regular expert testimony and source on public opinion for journalists, elected officials, and NGO leadership\"\n$cursor = Add-ParagraphAfter $cursor \"\u2022 Redistricting analysis used in court cases with rigorous methodology and expert testimony\"\n"}
